$wb = $excel.ActiveWorkbook
$refSheet = $wb.Worksheets.Item("References")
$newSheet = $wb.Worksheets.Add($refSheet)
$newSheet.Name = "Node.js"
Write-Output $wb.Worksheets.Count
for ($i=1; $i -le $wb.Worksheets.Count; $i++) {
    Write-Output "$i : $($wb.Worksheets.Item($i).Name)"
}
